$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename "美术（书法）教师" -> "美术教师" (cell A15 keeps its position in the list)
$ws.Range("A15").Value = "美术教师"

# Update the saved selection to A15 (was A19)
$ws.Range("A15").Select()
